$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (column D) values, preserving text storage
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "309.54"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "36.58"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.098"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07691"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.397"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.325"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.851"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.948"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9235"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1137"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1878"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.08759"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03361"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09533"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001382"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.005980"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.358"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3433"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.261"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1291"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04343"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001203"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004250"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001331"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0002903"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02123"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05016"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007491"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1348"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.008421"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002071"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.007718"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006360"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002868"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001691"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002102"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002002"

# Update Volume(1h) (column E) values, preserving text storage
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-1.11%"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-2.82%"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.84%"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-2.83%"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-0.55%"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.65%"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-3.58%"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.91%"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.20%"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-7.76%"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-2.11%"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-4.31%"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "1.43%"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.67%"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.16%"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "4.97%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-4.63%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.35%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "19.16%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1.44%"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-10.79%"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.55%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-3.71%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-1.18%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "9.06%"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-4.93%"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-2.17%"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "0.71%"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-0.93%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-3.60%"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "3.02%"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-10.53%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-5.52%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.07%"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-14.42%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "40.92%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.07%"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.07%"
